$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 117 (weekly update). Rows 117-139 shift down to 118-140,
# and the sheet's used range grows from A1:R139 to A1:R140.
$ws.Rows.Item(117).Insert()

# Populate the newly inserted row 117 with this week's data.
$ws.Cells.Item(117, 1).Value = 8
$ws.Cells.Item(117, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(117, 3).Value = "Coquimbo"
$ws.Cells.Item(117, 4).Value = "2021-12-13"
$ws.Cells.Item(117, 5).Value = 4
$ws.Cells.Item(117, 6).Value = 100112037
$ws.Cells.Item(117, 7).Value = "Cebollín"
$ws.Cells.Item(117, 8).Value = "Sin especificar"
$ws.Cells.Item(117, 9).Value = "Primera"
$ws.Cells.Item(117, 10).Value = 3000
$ws.Cells.Item(117, 11).Value = 900
$ws.Cells.Item(117, 12).Value = 1000
$ws.Cells.Item(117, 13).Value = 950
$ws.Cells.Item(117, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(117, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(117, 16).Value = 158
$ws.Cells.Item(117, 17).Value = 6
$ws.Cells.Item(117, 18).Value = "Hortaliza"
